{"js": "// Remove the unwanted trailing image: the last inline picture in the\n// document, which sits by itself in the final paragraph. The paragraph\n// mark itself (and any bookmark anchored to it, e.g. \"_GoBack\") must be\n// left in place - only the picture run goes away.\n\nconst doc = context.document;\nconst body = doc.body;\n\nconst pictures = body.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\nif (pictures.items.length > 0) {\n  const lastPicture = pictures.items[pictures.items.length - 1];\n  const paragraph = lastPicture.paragraph;\n  const paragraphRange = paragraph.getRange();\n\n  // Bookmarks whose name starts with \"_\" (like Word's own \"_GoBack\")\n  // are hidden and need includeHidden=true to be discovered.\n  const bookmarksResult = paragraphRange.getBookmarks(true);\n  await context.sync();\n\n  const bookmarkNames = bookmarksResult.value || [];\n\n  // Work around this paragraph's bookmark(s) otherwise getting dropped\n  // when the picture run is deleted: pull them out first, delete the\n  // picture, then restore the bookmark(s) onto the now-empty paragraph.\n  for (const name of bookmarkNames) {\n    doc.deleteBookmark(name);\n  }\n  await context.sync();\n\n  lastPicture.delete();\n  await context.sync();\n\n  if (bookmarkNames.length > 0) {\n    const restoreRange = paragraph.getRange();\n    for (const name of bookmarkNames) {\n      restoreRange.insertBookmark(name);\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Remove the unwanted trailing image: the last inline picture in the\n# document, which sits by itself in the final paragraph. The paragraph\n# mark itself (and any bookmark anchored to it, e.g. Word's own hidden\n# \"_GoBack\" bookmark) must be left in place - only the picture goes away.\n\n$d = $word.ActiveDocument\n\n$shapeCount = $d.InlineShapes.Count\nif ($shapeCount -gt 0) {\n    $lastShape = $d.InlineShapes.Item($shapeCount)\n    $shapeRange = $lastShape.Range\n    $paragraph = $shapeRange.Paragraphs.Item(1)\n\n    # Collect bookmark names anchored in this paragraph. Regular\n    # bookmarks show up in the Bookmarks collection, but Word hides\n    # bookmarks whose name starts with \"_\" (such as its own \"_GoBack\")\n    # from enumeration/Count, so also explicitly probe for that one.\n    $bookmarkNames = @()\n\n    $paraRange = $paragraph.Range\n    for ($i = 1; $i -le $paraRange.Bookmarks.Count; $i++) {\n        $bookmarkNames += $paraRange.Bookmarks.Item($i).Name\n    }\n\n    if ($d.Bookmarks.Exists(\"_GoBack\")) {\n        $goBackRange = $d.Bookmarks.Item(\"_GoBack\").Range\n        if ($goBackRange.Start -ge $paragraph.Range.Start -and $goBackRange.End -le $paragraph.Range.End) {\n            if (-not ($bookmarkNames -contains \"_GoBack\")) {\n                $bookmarkNames += \"_GoBack\"\n            }\n        }\n    }\n\n    foreach ($name in $bookmarkNames) {\n        if ($d.Bookmarks.Exists($name)) {\n            $d.Bookmarks.Item($name).Delete()\n        }\n    }\n\n    # Now remove the picture itself; the paragraph mark stays behind.\n    $lastShape.Delete()\n\n    # Restore the bookmark(s) onto the now-empty paragraph so they keep\n    # marking the same location.\n    if ($bookmarkNames.Count -gt 0) {\n        $restoreParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n        foreach ($name in $bookmarkNames) {\n            $d.Bookmarks.Add($name, $restoreParagraph.Range)\n        }\n    }\n}\n"}
